$d = $word.ActiveDocument

# Resume formatting: spacing around section headers.
#
# Section header paragraphs ("SUMMARY", "PRIMARY SKILLS",
# "PROFESSIONAL EXPERIENCE", "EDUCATION") currently carry only
# SpaceBefore = 6pt (120 twips). They all gain an explicit
# SpaceAfter = 4pt (80 twips); SUMMARY (the very first header in the
# document) has its SpaceBefore bumped to 8pt (160 twips) while the
# remaining headers go to 10pt (200 twips).
#
# All other non-bullet paragraphs (job-title lines, the summary blurb,
# and the education entries) currently carry SpaceAfter = 2pt
# (40 twips); those get bumped to 3pt (60 twips). Bulleted paragraphs
# (SpaceAfter = 1pt / 20 twips) are left untouched.

$isFirstHeader = $true

foreach ($p in $d.Paragraphs) {
    $before = $p.Format.SpaceBefore
    $after = $p.Format.SpaceAfter

    if ($before -eq 6) {
        # Section header.
        if ($isFirstHeader) {
            $p.Format.SpaceBefore = 8
            $isFirstHeader = $false
        } else {
            $p.Format.SpaceBefore = 10
        }
        $p.Format.SpaceAfter = 4
    } elseif ($after -eq 2) {
        # Job-title / summary-body / education paragraph.
        $p.Format.SpaceAfter = 3
    }
}
